$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D is treated as text so numeric-looking values (e.g. "0.999")
# are not auto-converted to numbers by Excel, keeping them as text like the source data.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$updates = @{
    2 = @{ 'D'='58.032.77'; 'E'='  -0.43%  ' }
    3 = @{ 'D'='2.453.83'; 'E'='  -1.31%  ' }
    4 = @{ 'E'='  -0.40%  ' }
    5 = @{ 'D'='525.47'; 'E'='  +0.72%  ' }
    6 = @{ 'D'='131.98'; 'E'='  -0.98%  ' }
    7 = @{ 'D'='0.999' }
    8 = @{ 'D'='0.565'; 'E'='  +0.53%  ' }
    9 = @{ 'D'='2.457.63'; 'E'='  -2.62%  ' }
    10 = @{ 'D'='0.0983'; 'E'='  +0.72%  ' }
    11 = @{ 'D'='0.152' }
    12 = @{ 'D'='5.01'; 'E'='  -3.09%  ' }
    13 = @{ 'E'='  -2.29%  ' }
    14 = @{ 'D'='2.887.38'; 'E'='  -2.80%  ' }
    15 = @{ 'D'='57.938.86'; 'E'='  -0.69%  ' }
    16 = @{ 'D'='21.80'; 'E'='  -1.47%  ' }
    17 = @{ 'E'='  -1.07%  ' }
    18 = @{ 'D'='2.454.94'; 'E'='  -3.01%  ' }
    19 = @{ 'D'='10.35'; 'E'='  -2.95%  ' }
    20 = @{ 'D'='4.13'; 'E'='  -0.72%  ' }
    21 = @{ 'D'='313.12'; 'E'='  -2.75%  ' }
    22 = @{ 'D'='6.11'; 'E'='  -1.13%  ' }
    23 = @{ 'E'='  -0.07%  ' }
    24 = @{ 'D'='65.00'; 'E'='  +0.54%  ' }
    25 = @{ 'D'='0.404'; 'E'='  -0.67%  ' }
    26 = @{ 'D'='0.999'; 'E'='  -0.04%  ' }
    27 = @{ 'E'='  -1.59%  ' }
    28 = @{ 'D'='7.25'; 'E'='  -2.16%  ' }
    29 = @{ 'D'='173.57'; 'E'='  +3.03%  ' }
    30 = @{ 'D'='0.0₃0740'; 'E'='  -1.67%  ' }
    31 = @{ 'D'='1.70'; 'E'='  -0.99%  ' }
    32 = @{ 'D'='6.26'; 'E'='  -0.59%  ' }
    33 = @{ 'E'='  -3.76%  ' }
    34 = @{ 'D'='0.998'; 'E'='  -0.01%  ' }
    35 = @{ 'D'='0.997'; 'E'='  +0.37%  ' }
    36 = @{ 'D'='17.85'; 'E'='  -1.90%  ' }
    37 = @{ 'D'='1.20'; 'E'='  -4.43%  ' }
    38 = @{ 'D'='3.81'; 'E'='  -2.80%  ' }
    39 = @{ 'B'='SuiNetwork'; 'C'='https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'; 'D'='0.816'; 'E'='  +5.99%  ' }
    40 = @{ 'B'='OKB'; 'C'='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; 'D'='36.28'; 'E'='  -0.33%  ' }
    41 = @{ 'B'='Stacks'; 'C'='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; 'D'='1.46'; 'E'='  -1.91%  ' }
    42 = @{ 'D'='3.42'; 'E'='  -0.78%  ' }
    43 = @{ 'B'='Bittensor'; 'C'='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; 'D'='262.93'; 'E'='  -5.08%  ' }
    44 = @{ 'B'='Mantle'; 'C'='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; 'D'='0.588'; 'E'='  -1.78%  ' }
    45 = @{ 'E'='  -3.59%  ' }
    46 = @{ 'D'='0.0922'; 'E'='  +0.33%  ' }
    47 = @{ 'D'='122.42'; 'E'='  -6.13%  ' }
    48 = @{ 'D'='0.0497'; 'E'='  -0.64%  ' }
    49 = @{ 'D'='0.0212'; 'E'='  -0.57%  ' }
    50 = @{ 'D'='16.97'; 'E'='  -4.16%  ' }
    51 = @{ 'D'='16.32'; 'E'='  -3.51%  ' }
}

foreach ($rowKey in $updates.Keys) {
    $rowNum = [int]$rowKey
    $cellVals = $updates[$rowKey]
    foreach ($col in $cellVals.Keys) {
        $ws.Range("$col$rowNum").Value = $cellVals[$col]
    }
}

# Reset style on the price column so no lingering number-format style is left
# attached to cells (keeps cells using the default/normal style like the rest of the sheet).
$priceRange.Style = "Normal"
